# Auto-generated: update Sheets via scheduled runner
# Refreshes market-price derived columns (H:N) for affected leve rows
# across all eight job sheets. Values only - no formulas in this workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    98 = @{ H=533999.7; I=589766.5; J=4215; K=589766.5; L=4215; M=-588268.5; N=-7211 }
    122 = @{ H=533999.7; I=589766.5; J=4215; K=1769299.5; L=12645; M=-1766849.5; N=-17545 }
    137 = @{ H=33334804; I=47620052; J=2564.889; K=142860156; L=7694.667; M=-142857606; N=-12794.667 }
    138 = @{ H=4632305; I=2472443.5; J=5130735; K=7417330.5; L=15392205; M=-7412190.5; N=-15402485 }
    141 = @{ H=2786.389; I=2362.0588; J=10000; K=7086.176399999999; L=30000; M=-1906.176399999999; N=-40360 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    2 = @{ H=2082.5715; I=1922.2667; J=2483.3333; K=1922.2667; L=2483.3333; M=-1809.2667; N=-2709.3333 }
    32 = @{ H=3505.985; I=2031.9818; J=10261.833; K=2031.9818; L=10261.833; M=-1744.9818; N=-10835.833 }
    45 = @{ H=1614; I=1302.4; J=2003.5; K=1302.4; L=2003.5; M=-925.4000000000001; N=-2757.5 }
    61 = @{ H=2889.394; I=1760.8889; J=4243.6; K=1760.8889; L=4243.6; M=-1548.8889; N=-4667.6 }
    102 = @{ H=883.1667; I=859.8; K=859.8; M=762.2 }
    116 = @{ H=2082.5715; I=1922.2667; J=2483.3333; K=1922.2667; L=2483.3333; M=371.7333000000001; N=-7071.3333 }
    122 = @{ H=2036.3636; I=2250; J=1988.8889; K=6750; L=5966.6667; M=-4300; N=-10866.6667 }
    123 = @{ H=34071; J=34071; L=34071; N=-43871 }
    124 = @{ H=18666.666; J=18666.666; L=18666.666; N=-28486.666 }
    125 = @{ H=33032.855; J=33032.855; L=33032.855; N=-42872.855 }
    132 = @{ H=2427.182; I=2015.9259; K=6047.7777; M=-3517.7777 }
    136 = @{ H=2889.394; I=1760.8889; J=4243.6; K=5282.6667; L=12730.8; M=-2732.6667; N=-17830.8 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    3 = @{ H=2082.5715; I=1922.2667; J=2483.3333; K=1922.2667; L=2483.3333; M=-1808.2667; N=-2711.3333 }
    107 = @{ H=1176.5555; I=1263; J=1068.5; K=1263; L=1068.5; M=657; N=-4908.5 }
    109 = @{ H=30000; J=30000; L=30000; N=-32774 }
    124 = @{ H=43660; J=43660; L=43660; N=-53480 }
    134 = @{ H=2383.9512; I=1596.7354; J=6207.5713; K=4790.206200000001; L=18622.7139; M=-2255.206200000001; N=-23692.7139 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    22 = @{ H=240.94444; I=197.8125; J=586; K=197.8125; L=586; M=152.1875; N=-1286 }
    31 = @{ H=1494.0834; I=1198; J=2382.3333; K=1198; L=2382.3333; M=-903; N=-2972.3333 }
    34 = @{ H=1494.0834; I=1198; J=2382.3333; K=1198; L=2382.3333; M=-996; N=-2786.3333 }
    58 = @{ H=1797.7241; I=1209.9048; K=1209.9048; M=-1006.9048 }
    122 = @{ H=1922; I=1230; J=3421.3333; K=3690; L=10263.9999; M=-1240; N=-15163.9999 }
    136 = @{ H=1797.7241; I=1209.9048; K=3629.7144; M=-1079.7144 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    106 = @{ H=4200; J=4200; L=12600; N=-14492 }
    107 = @{ H=391; I=371.58334; J=414.3; K=1114.75002; L=1242.9; M=805.2499800000001; N=-5082.9 }
    113 = @{ H=12821273; I=469.55554; K=1408.66662; M=761.33338 }
    137 = @{ H=7219255.5; J=130506.875; L=391520.625; N=-401720.625 }
    138 = @{ H=800.4666999999999; I=800.4666999999999; K=2401.4001; M=2738.5999 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    122 = @{ H=2223733.2; I=2779034.5; J=2528; K=8337103.5; L=7584; M=-8334653.5; N=-12484 }
    135 = @{ H=166681630; J=166681630; L=166681630; N=-166691770 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    93 = @{ H=1836; I=950; J=2722; K=950; L=2722; M=298; N=-5218 }
    122 = @{ H=3961.5386; I=3750; J=4000; K=11250; L=12000; M=-8800; N=-16900 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    107 = @{ H=455.73914; I=508.25; K=1524.75; M=395.25 }
    122 = @{ H=2828.2666; I=2122.4; K=6367.200000000001; M=-3917.200000000001 }
    126 = @{ H=60682.117; I=73278.64; J=1898.3334; K=219835.92; L=5695.0002; M=-217365.92; N=-10635.0002 }
}
foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
